$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update parser-table rows (grammar productions / FIRST / FOLLOW sets).
# Rows 2-25 get their text refreshed (some FIRST/FOLLOW sets changed because
# of the new <repeat_statement> / Check..Repeat loop construct with boolean
# conditions); rows 26-34 are brand new rows documenting the new grammar
# rules.
$ws.Cells.Item(2, 1).Value = "<program> -> Begin Colon EndOfLine <declaration_section> <function_definition> End SemiColon | ε"
$ws.Cells.Item(2, 2).Value = "Begin"
$ws.Cells.Item(2, 3).Value = "EOF"
$ws.Cells.Item(3, 1).Value = "<type> -> Int | Double | String | Boolean"
$ws.Cells.Item(3, 2).Value = "Int, Double, String, Boolean"
$ws.Cells.Item(3, 3).Value = "Identifier"
$ws.Cells.Item(4, 1).Value = "<declaration_section> -> Declaration Colon EndOfLine <var_declaration> End SemiColon EndOfLine"
$ws.Cells.Item(4, 2).Value = "Declaration"
$ws.Cells.Item(4, 3).Value = "Int, Double, String, Boolean"
$ws.Cells.Item(5, 1).Value = "<var_declaration> -> <var_declaration_line> <var_declaration> | ε"
$ws.Cells.Item(5, 2).Value = "Int, Double, String, Boolean, ε"
$ws.Cells.Item(5, 3).Value = "End"
$ws.Cells.Item(6, 1).Value = "<var_declaration_line> -> <type> <identifier_list> EndOfLine"
$ws.Cells.Item(6, 2).Value = "Int, Double, String, Boolean"
$ws.Cells.Item(6, 3).Value = "Int, Double, String, Boolean, End"
$ws.Cells.Item(7, 1).Value = "<identifier_list> -> Identifier <identifier_tail>"
$ws.Cells.Item(7, 2).Value = "Identifier"
$ws.Cells.Item(7, 3).Value = "EndOfLine, SemiColon"
$ws.Cells.Item(8, 1).Value = "<identifier_tail> -> Comma Identifier <identifier_tail> | ε"
$ws.Cells.Item(8, 2).Value = "Comma, ε"
$ws.Cells.Item(8, 3).Value = "EndOfLine, SemiColon"
$ws.Cells.Item(9, 1).Value = "<function_definition> -> <type> Identifier LeftParen <parameter_list> RightParen Colon EndOfLine <statement_list> End SemiColon EndOfLine"
$ws.Cells.Item(9, 2).Value = "Int, Double, String, Boolean"
$ws.Cells.Item(9, 3).Value = "End"
$ws.Cells.Item(10, 1).Value = "<parameter_list> -> <parameter> <parameter_tail> | ε"
$ws.Cells.Item(10, 2).Value = "Int, Double, String, Boolean, ε"
$ws.Cells.Item(10, 3).Value = "RightParen"
$ws.Cells.Item(11, 1).Value = "<parameter> -> <type> Identifier"
$ws.Cells.Item(11, 2).Value = "Int, Double, String, Boolean"
$ws.Cells.Item(11, 3).Value = "Comma, RightParen"
$ws.Cells.Item(12, 1).Value = "<parameter_tail> -> Comma <parameter> <parameter_tail> | ε"
$ws.Cells.Item(12, 2).Value = "Comma, ε"
$ws.Cells.Item(12, 3).Value = "RightParen"
$ws.Cells.Item(13, 1).Value = "<statement_list> -> <statement> <statement_list> | ε"
$ws.Cells.Item(13, 2).Value = "Identifier, Input, Output, Return, Check, ε"
$ws.Cells.Item(13, 3).Value = "End"
$ws.Cells.Item(14, 1).Value = "<statement> -> <assignment_statement> | <input_statement> | <output_statement> | <return_statement> | <repeat_statement>"
$ws.Cells.Item(14, 2).Value = "Identifier, Input, Output, Return, Check"
$ws.Cells.Item(14, 3).Value = "Identifier, Input, Output, Return, Check, End"
$ws.Cells.Item(15, 1).Value = "<assignment_statement> -> Identifier Assignment <expression> EndOfLine"
$ws.Cells.Item(15, 2).Value = "Identifier"
$ws.Cells.Item(15, 3).Value = "Identifier, Input, Output, Return, Check, Break, Continue, End"
$ws.Cells.Item(16, 1).Value = "<expression> -> <mul_expre> <add_expre_tail>"
$ws.Cells.Item(16, 2).Value = "Identifier, IntLiteral, StringLiteral, DoubleLiteral, LeftParen"
$ws.Cells.Item(16, 3).Value = "EndOfLine, RightParen"
$ws.Cells.Item(17, 1).Value = "<add_expre_tail> -> Add <mul_expre> <add_expre_tail> | Subtract <mul_expre> <add_expre_tail> | ε"
$ws.Cells.Item(17, 2).Value = "Add, Subtract, ε"
$ws.Cells.Item(17, 3).Value = "EndOfLine, RightParen"
$ws.Cells.Item(18, 1).Value = "<mul_expre> -> <pow_expre> <mul_expre_tail>"
$ws.Cells.Item(18, 2).Value = "Identifier, IntLiteral, StringLiteral, DoubleLiteral, LeftParen"
$ws.Cells.Item(18, 3).Value = "Add, Subtract, EndOfLine, RightParen"
$ws.Cells.Item(19, 1).Value = "<mul_expre_tail> -> Multiply <pow_expre> <mul_expre_tail> | Divide <pow_expre> <mul_expre_tail> | Modulo <pow_expre> <mul_expre_tail> | ε"
$ws.Cells.Item(19, 2).Value = "Multiply, Divide, Modulo, ε"
$ws.Cells.Item(19, 3).Value = "Add, Subtract, EndOfLine, RightParen"
$ws.Cells.Item(20, 1).Value = "<pow_expre> -> <factor> Power <pow_expre> | <factor>"
$ws.Cells.Item(20, 2).Value = "Identifier, IntLiteral, StringLiteral, DoubleLiteral, LeftParen"
$ws.Cells.Item(20, 3).Value = "Multiply, Divide, Modulo, Add, Subtract, EndOfLine, RightParen"
$ws.Cells.Item(21, 1).Value = "<factor> -> Identifier | IntLiteral | StringLiteral | DoubleLiteral | LeftParen <expression> RightParen"
$ws.Cells.Item(21, 2).Value = "Identifier, IntLiteral, StringLiteral, DoubleLiteral, LeftParen"
$ws.Cells.Item(21, 3).Value = "Multiply, Divide, Modulo, Add, Subtract, EndOfLine, RightParen, Power"
$ws.Cells.Item(22, 1).Value = "<input_statement> -> Input Colon <identifier_list> SemiColon EndOfLine"
$ws.Cells.Item(22, 2).Value = "Input"
$ws.Cells.Item(22, 3).Value = "Identifier, Input, Output, Return, Check, Break, Continue, End"
$ws.Cells.Item(23, 1).Value = "<output_statement> -> Output Colon <output_target> SemiColon EndOfLine"
$ws.Cells.Item(23, 2).Value = "Output"
$ws.Cells.Item(23, 3).Value = "Identifier, Input, Output, Return, Check, Break, Continue, End"
$ws.Cells.Item(24, 1).Value = "<output_target> -> Identifier | StringLiteral | StringDoubleQuoteLiteral"
$ws.Cells.Item(24, 2).Value = "Identifier, StringLiteral, StringDoubleQuoteLiteral"
$ws.Cells.Item(24, 3).Value = "SemiColon"
$ws.Cells.Item(25, 1).Value = "<return_statement> -> Return <expression> EndOfLine"
$ws.Cells.Item(25, 2).Value = "Return"
$ws.Cells.Item(25, 3).Value = "Identifier, Input, Output, Return, Check, Break, Continue, End"
$ws.Cells.Item(26, 1).Value = "<repeat_statement> -> Check <condition> Colon EndOfLine <repeat_statement_list> Repeat SemiColon EndOfLine"
$ws.Cells.Item(26, 2).Value = "Check"
$ws.Cells.Item(26, 3).Value = "Identifier, Input, Output, Return, Check, Break, Continue, End"
$ws.Cells.Item(27, 1).Value = "<repeat_statement_list> -> <repeat_statement_line> <repeat_statement_list> | ε"
$ws.Cells.Item(27, 2).Value = "Identifier, Input, Output, Return, Check, Break, Continue, ε"
$ws.Cells.Item(27, 3).Value = "Repeat"
$ws.Cells.Item(28, 1).Value = "<repeat_statement_line> -> <assignment_statement> | <input_statement> | <output_statement> | <return_statement> | <repeat_statement> | Break EndOfLine | Continue EndOfLine"
$ws.Cells.Item(28, 2).Value = "Identifier, Input, Output, Return, Check, Break, Continue"
$ws.Cells.Item(28, 3).Value = "Identifier, Input, Output, Return, Check, Break, Continue, Repeat"
$ws.Cells.Item(29, 1).Value = "<condition> -> <bool_term> <bool_expr_tail>"
$ws.Cells.Item(29, 2).Value = "LogicalNot, True, False, Identifier, IntLiteral, DoubleLiteral, StringLiteral, LeftParen"
$ws.Cells.Item(29, 3).Value = "Colon"
$ws.Cells.Item(30, 1).Value = "<bool_expr_tail> -> LogicalAnd <bool_term> <bool_expr_tail> | LogicalOr <bool_term> <bool_expr_tail> | ε"
$ws.Cells.Item(30, 2).Value = "LogicalAnd, LogicalOr, ε"
$ws.Cells.Item(30, 3).Value = "Colon"
$ws.Cells.Item(31, 1).Value = "<bool_term> -> LogicalNot <basic_bool> | <basic_bool>"
$ws.Cells.Item(31, 2).Value = "LogicalNot, True, False, Identifier, IntLiteral, DoubleLiteral, StringLiteral, LeftParen"
$ws.Cells.Item(31, 3).Value = "LogicalAnd, LogicalOr, Colon"
$ws.Cells.Item(32, 1).Value = "<basic_bool> -> True | False | <operand> <compare_op> <operand> | LeftParen <bool_expr> RightParen"
$ws.Cells.Item(32, 2).Value = "True, False, Identifier, IntLiteral, DoubleLiteral, StringLiteral, LeftParen"
$ws.Cells.Item(32, 3).Value = "LogicalAnd, LogicalOr, Colon"
$ws.Cells.Item(33, 1).Value = "<operand> -> Identifier | IntLiteral | DoubleLiteral | StringLiteral | True | False"
$ws.Cells.Item(33, 2).Value = "Identifier, IntLiteral, DoubleLiteral, StringLiteral, True, False"
$ws.Cells.Item(33, 3).Value = "Equal, NotEqual, LessThan, GreaterThan, LessOrEqual, GreaterOrEqual, LogicalAnd, LogicalOr, Colon"
$ws.Cells.Item(34, 1).Value = "<compare_op> -> Equal | NotEqual | LessThan | GreaterThan | LessOrEqual | GreaterOrEqual"
$ws.Cells.Item(34, 2).Value = "Equal, NotEqual, LessThan, GreaterThan, LessOrEqual, GreaterOrEqual"
$ws.Cells.Item(34, 3).Value = "Identifier, IntLiteral, DoubleLiteral, StringLiteral, True, False"

# Column sizing / view state follow the widened text in columns A and B.
$ws.Columns.Item(1).ColumnWidth = 157
$ws.Columns.Item(2).ColumnWidth = 72.3

$excel.ActiveWindow.Zoom = 89
$ws.Range("C36").Select()
